$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B22 was stored as text "20"; convert it to a real number
$ws.Range("B22").Value = 20

# Append new row 23
$ws.Range("A23").Value = "Koemthay Tha"

# B23 should stay a text value "19" (not auto-converted to a number)
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "19"
$ws.Range("B23").Style = "Normal"

$ws.Range("C23").Value = "Male"
$ws.Range("D23").Value = "Kompong Cham"
$ws.Range("E23").Value = "Class B 2025"
$ws.Range("F23").Value = "image\4f092cd9e1354a2a91fc0256577f88b7.png"
